$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 193.88889
$ws.Range("J2").Value = 168.8
$ws.Range("L2").Value = 168.8
$ws.Range("N2").Value = -394.8
$ws.Range("H33").Value = 117.545456
$ws.Range("I33").Value = 119.75
$ws.Range("K33").Value = 119.75
$ws.Range("M33").Value = 109.25
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H64").Value = 3230
$ws.Range("J64").Value = 3287.5
$ws.Range("L64").Value = 3287.5
$ws.Range("N64").Value = -3783.5
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H67").Value = 3230
$ws.Range("J67").Value = 3287.5
$ws.Range("L67").Value = 3287.5
$ws.Range("N67").Value = -5003.5
$ws.Range("H76").Value = 1381872.9
$ws.Range("I76").Value = 2345181.5
$ws.Range("J76").Value = 5717.857
$ws.Range("K76").Value = 2345181.5
$ws.Range("L76").Value = 5717.857
$ws.Range("M76").Value = -2344866.5
$ws.Range("N76").Value = -6347.857
$ws.Range("H79").Value = 1381872.9
$ws.Range("I79").Value = 2345181.5
$ws.Range("J79").Value = 5717.857
$ws.Range("K79").Value = 2345181.5
$ws.Range("L79").Value = 5717.857
$ws.Range("M79").Value = -2344089.5
$ws.Range("N79").Value = -7901.857
$ws.Range("H86").Value = 1773.3334
$ws.Range("I86").Value = 1700
$ws.Range("J86").Value = 1846.6666
$ws.Range("K86").Value = 1700
$ws.Range("L86").Value = 1846.6666
$ws.Range("M86").Value = -577
$ws.Range("N86").Value = -4092.6666
$ws.Range("H89").Value = 1773.3334
$ws.Range("I89").Value = 1700
$ws.Range("J89").Value = 1846.6666
$ws.Range("K89").Value = 8500
$ws.Range("L89").Value = 9233.333000000001
$ws.Range("M89").Value = -2884
$ws.Range("N89").Value = -20465.333
$ws.Range("H138").Value = 3131.851
$ws.Range("J138").Value = 2905.037
$ws.Range("L138").Value = 8715.110999999999
$ws.Range("N138").Value = -18995.111

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5291.75
$ws.Range("J63").Value = 3843.3333
$ws.Range("L63").Value = 3843.3333
$ws.Range("N63").Value = -5215.3333
$ws.Range("H66").Value = 5291.75
$ws.Range("J66").Value = 3843.3333
$ws.Range("L66").Value = 19216.6665
$ws.Range("N66").Value = -26080.6665
$ws.Range("H74").Value = 1313.3077
$ws.Range("J74").Value = 3849.8333
$ws.Range("L74").Value = 3849.8333
$ws.Range("N74").Value = -5597.8333
$ws.Range("H77").Value = 1313.3077
$ws.Range("J77").Value = 3849.8333
$ws.Range("L77").Value = 19249.1665
$ws.Range("N77").Value = -27985.1665
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1592.4445
$ws.Range("I132").Value = 982.5
$ws.Range("K132").Value = 2947.5
$ws.Range("M132").Value = -417.5

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 107589.42
$ws.Range("I86").Value = 2207.6924
$ws.Range("K86").Value = 2207.6924
$ws.Range("M86").Value = -1084.6924
$ws.Range("H89").Value = 107589.42
$ws.Range("I89").Value = 2207.6924
$ws.Range("K89").Value = 11038.462
$ws.Range("M89").Value = -5422.462

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2590.9092
$ws.Range("I31").Value = 1757.7742
$ws.Range("J31").Value = 4577.615
$ws.Range("K31").Value = 1757.7742
$ws.Range("L31").Value = 4577.615
$ws.Range("M31").Value = -1462.7742
$ws.Range("N31").Value = -5167.615
$ws.Range("H34").Value = 2590.9092
$ws.Range("I34").Value = 1757.7742
$ws.Range("J34").Value = 4577.615
$ws.Range("K34").Value = 1757.7742
$ws.Range("L34").Value = 4577.615
$ws.Range("M34").Value = -1555.7742
$ws.Range("N34").Value = -4981.615
$ws.Range("H58").Value = 1632.2941
$ws.Range("I58").Value = 1388.4
$ws.Range("J58").Value = 1980.7142
$ws.Range("K58").Value = 1388.4
$ws.Range("L58").Value = 1980.7142
$ws.Range("M58").Value = -1185.4
$ws.Range("N58").Value = -2386.7142
$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678
$ws.Range("H136").Value = 1632.2941
$ws.Range("I136").Value = 1388.4
$ws.Range("J136").Value = 1980.7142
$ws.Range("K136").Value = 4165.200000000001
$ws.Range("L136").Value = 5942.142599999999
$ws.Range("M136").Value = -1615.200000000001
$ws.Range("N136").Value = -11042.1426

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10623.375
$ws.Range("I87").Value = 5831.1665
$ws.Range("K87").Value = 17493.4995
$ws.Range("M87").Value = -16245.4995
$ws.Range("H90").Value = 10623.375
$ws.Range("I90").Value = 5831.1665
$ws.Range("K90").Value = 52480.4985
$ws.Range("M90").Value = -46240.4985
$ws.Range("H131").Value = 776.16
$ws.Range("J131").Value = 805.337
$ws.Range("L131").Value = 2416.011
$ws.Range("N131").Value = -12496.011

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10636.667
$ws.Range("I70").Value = 15291.571
$ws.Range("J70").Value = 4119.8
$ws.Range("K70").Value = 15291.571
$ws.Range("L70").Value = 4119.8
$ws.Range("M70").Value = -15021.571
$ws.Range("N70").Value = -4659.8
$ws.Range("H73").Value = 10636.667
$ws.Range("I73").Value = 15291.571
$ws.Range("J73").Value = 4119.8
$ws.Range("K73").Value = 15291.571
$ws.Range("L73").Value = 4119.8
$ws.Range("M73").Value = -14355.571
$ws.Range("N73").Value = -5991.8
$ws.Range("H80").Value = 3023.6667
$ws.Range("J80").Value = 2916.6667
$ws.Range("L80").Value = 2916.6667
$ws.Range("N80").Value = -4912.6667
$ws.Range("H83").Value = 3023.6667
$ws.Range("J83").Value = 2916.6667
$ws.Range("L83").Value = 14583.3335
$ws.Range("N83").Value = -24567.3335
$ws.Range("H97").Value = 867.7273
$ws.Range("J97").Value = 953
$ws.Range("L97").Value = 953
$ws.Range("N97").Value = -1945
$ws.Range("H102").Value = 3375.3
$ws.Range("I102").Value = 3594
$ws.Range("J102").Value = 2500.5
$ws.Range("K102").Value = 3594
$ws.Range("L102").Value = 2500.5
$ws.Range("M102").Value = -1972
$ws.Range("N102").Value = -5744.5
$ws.Range("H122").Value = 2156.4546
$ws.Range("I122").Value = 1678.5
$ws.Range("K122").Value = 5035.5
$ws.Range("M122").Value = -2585.5

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5737.3335
$ws.Range("I32").Value = 5737.3335
$ws.Range("K32").Value = 5737.3335
$ws.Range("M32").Value = -5420.3335
$ws.Range("H40").Value = 9428
$ws.Range("J40").Value = 14999.75
$ws.Range("L40").Value = 14999.75
$ws.Range("N40").Value = -15271.75
$ws.Range("H122").Value = 8833.916999999999
$ws.Range("I122").Value = 7600.7
$ws.Range("K122").Value = 22802.1
$ws.Range("M122").Value = -20352.1
$ws.Range("H136").Value = 3373.5454
$ws.Range("I136").Value = 1674.5454
$ws.Range("J136").Value = 5072.5454
$ws.Range("K136").Value = 5023.6362
$ws.Range("L136").Value = 15217.6362
$ws.Range("M136").Value = -2473.6362
$ws.Range("N136").Value = -20317.6362

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1431.3636
$ws.Range("I81").Value = 1158.125
$ws.Range("K81").Value = 2316.25
$ws.Range("M81").Value = -1255.25
$ws.Range("H84").Value = 1431.3636
$ws.Range("I84").Value = 1158.125
$ws.Range("K84").Value = 11581.25
$ws.Range("M84").Value = -6277.25
$ws.Range("H100").Value = 626.63635
$ws.Range("I100").Value = 421.16666
$ws.Range("K100").Value = 842.33332
$ws.Range("M100").Value = -301.33332
$ws.Range("H105").Value = 31485.2
$ws.Range("J105").Value = 31485.2
$ws.Range("L105").Value = 31485.2
$ws.Range("N105").Value = -38473.2
$ws.Range("H108").Value = 55124.5
$ws.Range("J108").Value = 55124.5
$ws.Range("L108").Value = 55124.5
$ws.Range("N108").Value = -62804.5
$ws.Range("H136").Value = 2625.7666
$ws.Range("I136").Value = 2747.75
$ws.Range("J136").Value = 2486.3572
$ws.Range("K136").Value = 8243.25
$ws.Range("L136").Value = 7459.071599999999
$ws.Range("M136").Value = -5693.25
$ws.Range("N136").Value = -12559.0716
